$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two extra rows (20 and 21) first so dimension shrinks to A1:E19
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(20).Delete()

# Update data rows 2..19 with the new recommendation values
$ws.Cells.Item(2, 1).Value = 419
$ws.Cells.Item(2, 2).Value = 904
$ws.Cells.Item(2, 3).Value = 3.579442714350294
$ws.Cells.Item(2, 4).Value = 4.40346779511421
$ws.Cells.Item(2, 5).Value = "Rear Window (1954)"

$ws.Cells.Item(3, 1).Value = 419
$ws.Cells.Item(3, 2).Value = 1204
$ws.Cells.Item(3, 3).Value = 3.579442714350294
$ws.Cells.Item(3, 4).Value = 4.383534274137288
$ws.Cells.Item(3, 5).Value = "Lawrence of Arabia (1962)"

$ws.Cells.Item(4, 1).Value = 419
$ws.Cells.Item(4, 2).Value = 246
$ws.Cells.Item(4, 3).Value = 3.579442714350294
$ws.Cells.Item(4, 4).Value = 4.381324393346412
$ws.Cells.Item(4, 5).Value = "Hoop Dreams (1994)"

$ws.Cells.Item(5, 1).Value = 419
$ws.Cells.Item(5, 2).Value = 7361
$ws.Cells.Item(5, 3).Value = 3.579442714350294
$ws.Cells.Item(5, 4).Value = 4.378079465621082
$ws.Cells.Item(5, 5).Value = "Eternal Sunshine of the Spotless Mind (2004)"

$ws.Cells.Item(6, 1).Value = 419
$ws.Cells.Item(6, 2).Value = 2324
$ws.Cells.Item(6, 3).Value = 3.579442714350294
$ws.Cells.Item(6, 4).Value = 4.37424505401661
$ws.Cells.Item(6, 5).Value = "Life Is Beautiful (La Vita è bella) (1997)"

$ws.Cells.Item(7, 1).Value = 256
$ws.Cells.Item(7, 2).Value = 7371
$ws.Cells.Item(7, 3).Value = 3.579442714350294
$ws.Cells.Item(7, 4).Value = 4.519098607437904
$ws.Cells.Item(7, 5).Value = "Dogville (2003)"

$ws.Cells.Item(8, 1).Value = 256
$ws.Cells.Item(8, 2).Value = 6016
$ws.Cells.Item(8, 3).Value = 3.579442714350294
$ws.Cells.Item(8, 4).Value = 4.505543143517541
$ws.Cells.Item(8, 5).Value = "City of God (Cidade de Deus) (2002)"

$ws.Cells.Item(9, 1).Value = 256
$ws.Cells.Item(9, 2).Value = 2542
$ws.Cells.Item(9, 3).Value = 3.579442714350294
$ws.Cells.Item(9, 4).Value = 4.478808819824815
$ws.Cells.Item(9, 5).Value = "Lock, Stock & Two Smoking Barrels (1998)"

$ws.Cells.Item(10, 1).Value = 256
$ws.Cells.Item(10, 2).Value = 1242
$ws.Cells.Item(10, 3).Value = 3.579442714350294
$ws.Cells.Item(10, 4).Value = 4.474695418345246
$ws.Cells.Item(10, 5).Value = "Glory (1989)"

$ws.Cells.Item(11, 1).Value = 256
$ws.Cells.Item(11, 2).Value = 48516
$ws.Cells.Item(11, 3).Value = 3.579442714350294
$ws.Cells.Item(11, 4).Value = 4.462593355387727
$ws.Cells.Item(11, 5).Value = "Departed, The (2006)"

$ws.Cells.Item(12, 1).Value = 70
$ws.Cells.Item(12, 2).Value = 898
$ws.Cells.Item(12, 3).Value = 3.579442714350294
$ws.Cells.Item(12, 4).Value = 4.653532122802486
$ws.Cells.Item(12, 5).Value = "Philadelphia Story, The (1940)"

$ws.Cells.Item(13, 1).Value = 70
$ws.Cells.Item(13, 2).Value = 750
$ws.Cells.Item(13, 3).Value = 3.579442714350294
$ws.Cells.Item(13, 4).Value = 4.649597859379635
$ws.Cells.Item(13, 5).Value = "Dr. Strangelove or: How I Learned to Stop Worrying and Love the Bomb (1964)"

$ws.Cells.Item(14, 1).Value = 70
$ws.Cells.Item(14, 2).Value = 1136
$ws.Cells.Item(14, 3).Value = 3.579442714350294
$ws.Cells.Item(14, 4).Value = 4.634275890776808
$ws.Cells.Item(14, 5).Value = "Monty Python and the Holy Grail (1975)"

$ws.Cells.Item(15, 1).Value = 70
$ws.Cells.Item(15, 2).Value = 2959
$ws.Cells.Item(15, 3).Value = 3.579442714350294
$ws.Cells.Item(15, 4).Value = 4.624144333833454
$ws.Cells.Item(15, 5).Value = "Fight Club (1999)"

$ws.Cells.Item(16, 1).Value = 345
$ws.Cells.Item(16, 2).Value = 318
$ws.Cells.Item(16, 3).Value = 3.579442714350294
$ws.Cells.Item(16, 4).Value = 4.562637723843174
$ws.Cells.Item(16, 5).Value = "Shawshank Redemption, The (1994)"

$ws.Cells.Item(17, 1).Value = 345
$ws.Cells.Item(17, 2).Value = 858
$ws.Cells.Item(17, 3).Value = 3.579442714350294
$ws.Cells.Item(17, 4).Value = 4.417489923360288
$ws.Cells.Item(17, 5).Value = "Godfather, The (1972)"

$ws.Cells.Item(18, 1).Value = 241
$ws.Cells.Item(18, 2).Value = 1104
$ws.Cells.Item(18, 3).Value = 3.579442714350294
$ws.Cells.Item(18, 4).Value = 4.20735890635388
$ws.Cells.Item(18, 5).Value = "Streetcar Named Desire, A (1951)"

$ws.Cells.Item(19, 1).Value = 241
$ws.Cells.Item(19, 2).Value = 1208
$ws.Cells.Item(19, 3).Value = 3.579442714350294
$ws.Cells.Item(19, 4).Value = 4.204941180277887
$ws.Cells.Item(19, 5).Value = "Apocalypse Now (1979)"
